$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

$ws.Range("A11").Value = "replace"
$ws.Range("B11").Value = "Replace basic"

$ws.Range("A12").Value = "search"
$ws.Range("B12").Value = "Search basic"

$ws.Range("C12").Value = "/keyword     //serach a keyword`nn               //search next hit`nN              //search last hit"
$ws.Range("C11").Value = ":s/foo/bar/g          //do substitile in current line`n:%s/foo/bar/g       //do substitile in all lines"

$ws.Range("C11:C12").WrapText = $true
$ws.Rows.Item(11).RowHeight = 30
$ws.Rows.Item(12).RowHeight = 45

$ws.Range("C12").Select()
